$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").ClearContents()
$ws.Range("B1").ClearContents()
$ws.Rows.Item(1).RowHeight = 18.75
